$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing ID values (201234 -> 403121)
$ws.Range("A1").Value = 403121
$ws.Range("F2").Value = 403121

# Row 3: Likhitha Pulluru
$ws.Range("A3").Value = 501302
$ws.Range("B3").Value = "pullurul"
$ws.Range("C3").Value = "Likhitha Pulluru"
$ws.Range("D3").Value = "Likhitha.Pulluru@cdk.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:Likhitha.Pulluru@cdk.com") | Out-Null
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = 403121

# Row 4: Sowmya Golla
$ws.Range("A4").Value = 501278
$ws.Range("C4").Value = "Sowmya Golla"
$ws.Range("D4").Value = "Sowmya.Golla@cdk.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:Sowmya.Golla@cdk.com") | Out-Null
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("B4").Value = "gollas"
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 403121

$ws.Range("B4").Select() | Out-Null

# Page setup (paper size / orientation), as seen after interacting with Page Setup dialog
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
